$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Same domainEx")

$ws.Range("G4").Value = "Distance"

$ws.Range("C5").Value = 23
$ws.Range("D5").Value = 9
$ws.Range("E5").Value = 1
$ws.Range("F5").Value = 53
$ws.Range("G5").Value = 3.28

$ws.Range("C6").Value = 18
$ws.Range("D6").Value = 16
$ws.Range("E6").Value = 8

$ws.Range("G6").Select()
